$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.847.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "'2.403.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'570.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").Value = "'139.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.22%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "'2.384.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.14%  "
$ws.Range("D10").Value = "'0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "'5.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").Value = "'25.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.0000170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "'2.818.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("D17").Value = "'60.755.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "'2.397.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "'321.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'4.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'1.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.03%  "
$ws.Range("D26").Value = "'64.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "'8.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.57%  "
$ws.Range("D28").Value = "'577.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.24%  "
$ws.Range("D30").Value = "'0.0₃0906"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.31%  "
$ws.Range("D31").Value = "'7.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  -7.05%  "
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("E34").Value = "  -8.11%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").Value = "'4.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.95%  "
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("D39").Value = "'147.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").Value = "'18.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").Value = "'5.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.32%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'41.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("E44").Value = "  -5.35%  "
$ws.Range("E45").Value = "  -6.03%  "
$ws.Range("D46").Value = "'0.0₆0287"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +20.81%  "
$ws.Range("D47").Value = "'140.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("E48").Value = "  -4.61%  "
$ws.Range("D49").Value = "'0.585"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("D50").Value = "'0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").Value = "'19.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.39%  "
